$wb = $excel.ActiveWorkbook

# --- Teilnehmer 8 (Alex) --- copied from "Teilnehmer 7 (Johannes)"
$src = $wb.Worksheets.Item("Teilnehmer 7 (Johannes)")
$src.Copy($null, $src)
$alex = $wb.Worksheets.Item($src.Index + 1)
$alex.Name = "Teilnehmer 8 (Alex)"

$alex.Range("B6").Value = "Uhrzeit: 18:08"
$alex.Range("B28").Value = "Containerauswertungsreihenfolge wichtig für's Verständnis"
$alex.Range("B30").Value = "nutzt selected transformation statt visibile"
$alex.Range("B21").Value = "Uhrzeit: 18:51"

# --- Teilnehmer 9 (Janos) --- copied from "Teilnehmer 8 (Alex)"
$alex.Copy($null, $alex)
$janos = $wb.Worksheets.Item($alex.Index + 1)
$janos.Name = "Teilnehmer 9 (Janos)"

$janos.Range("B28").Value = ""
$janos.Range("B30").Value = ""

$janos.Range("B6").Value = "Uhrzeit: 13:20 Uhr"
$janos.Range("B31").Value = "Proband erwartet dass sich Transformationen überlagern (visible und transparent)"
$janos.Range("B33").Value = "Will gesamte PhoneUtils Klasse sehen und Beziehungen selected"
$janos.Range("B34").Value = "relationen aufteilen in innerhalb & außerhalb"
$janos.Range("B21").Value = "Uhrzeit: 14:16 Uhr"

# --- Teilnehmer 10 (Lisa) --- copied from "Teilnehmer 9 (Janos)"
$janos.Copy($null, $janos)
$lisa = $wb.Worksheets.Item($janos.Index + 1)
$lisa.Name = "Teilnehmer 10 (Lisa)"

$lisa.Range("B31").Value = ""
$lisa.Range("B33").Value = ""
$lisa.Range("B34").Value = ""

$lisa.Range("B6").Value = "Uhrzeit: 15:00 Uhr"
$lisa.Range("B29").Value = "vorschlagseinschränkung verwirrt"
$lisa.Range("B21").Value = "Uhrzeit: 15:45 Uhr"

$lisa.Select()
